$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the unit-number header cells: "100.x" -> "101.x"
$ws.Range("A3").Value = 101.1
$ws.Range("A36").Value = "101.5"
$ws.Range("A28").Value = "101.4"
$ws.Range("A21").Value = "101.3"
$ws.Range("A13").Value = "101.2"

# Move the active selection to A4 (as recorded in the saved view state)
$ws.Range("A4").Select()
